$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of fiscal-position / tax mapping data (AT, BE, NL), mirroring
# the existing EU / XX / RC / SP / LI / 4% blocks already on the sheet.
$data = @(
    @{ Row = 12; A = "z0bug.fiscalpos_at_1"; C = "z0bug.fiscalpos_at"; D = "z0bug.tax_22v";  E = "z0bug.tax_eu-1-AT-v" },
    @{ Row = 13; A = "z0bug.fiscalpos_at_2"; C = "z0bug.fiscalpos_at"; D = "z0bug.tax_10v";  E = "z0bug.tax_eu-3-AT-v" },
    @{ Row = 14; A = "z0bug.fiscalpos_be_1"; C = "z0bug.fiscalpos_be"; D = "z0bug.tax_22v";  E = "z0bug.tax_eu-1-BE-v" },
    @{ Row = 15; A = "z0bug.fiscalpos_be_2"; C = "z0bug.fiscalpos_be"; D = "z0bug.tax_10v";  E = "z0bug.tax_eu-3-BE-v" },
    @{ Row = 16; A = "z0bug.fiscalpos_nl_1"; C = "z0bug.fiscalpos_nl"; D = "z0bug.tax_22v";  E = "z0bug.tax_eu-1-NL-v" },
    @{ Row = 17; A = "z0bug.fiscalpos_nl_2"; C = "z0bug.fiscalpos_nl"; D = "z0bug.tax_10v";  E = "z0bug.tax_eu-3-NL-v" }
)

foreach ($r in $data) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

$ws.Range("E18").Select()
